# Apply data updates to Sheet1, then delete row 3 (which removes it from the
# used range / dimension, shrinking it from A1:BD3 to A1:BD2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update individual cell values on row 2
$ws.Range("G2").Value  = 1.73
$ws.Range("I2").Value  = 4.2
$ws.Range("K2").Value  = 2.6
$ws.Range("L2").Value  = 4
$ws.Range("M2").Value  = 1.01
$ws.Range("N2").Value  = 23
$ws.Range("W2").Value  = 13
$ws.Range("X2").Value  = 12
$ws.Range("AF2").Value = 29
$ws.Range("AI2").Value = 26
$ws.Range("BB2").Value = 51

# Remove row 3 entirely
$ws.Rows.Item(3).Delete()
